$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.230642914772034
$ws.Range("B1").Value = 2.50421404838562
$ws.Range("C1").Value = 4.441522598266602
$ws.Range("D1").Value = 2.523002624511719
$ws.Range("E1").Value = 1.076730489730835
